# Update "苏州-漫展信息" workbook to reflect the gh-pages data refresh
# (commit: "Update gh-pages to output generated at 456a3b4").
#
# Changes on both the "展览" sheet and the "全部类型" sheet (same events,
# the latter has one extra leading row so everything is shifted by +1):
#  - "想去人数" (column F) counts bumped for several events.
#  - The "苏州·国风宠物-cosplay展" event was cancelled: its name now carries
#    a "（取消）" suffix and its "最低票价" (column G) switched from a
#    numeric price to the text "不可售" (not for sale).

$wb = $excel.ActiveWorkbook

# ---- Sheet "展览" ----
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F4").Value = 94
$ws1.Range("F6").Value = 45
$ws1.Range("F7").Value = 2676

$ws1.Range("C8").Value = "苏州·国风宠物-cosplay展（取消）"
$ws1.Range("G8").Value = "不可售"

$ws1.Range("F9").Value = 259
$ws1.Range("F10").Value = 116
$ws1.Range("F11").Value = 10043
$ws1.Range("F13").Value = 256
$ws1.Range("F14").Value = 8
$ws1.Range("F15").Value = 620
$ws1.Range("F16").Value = 11750
$ws1.Range("F17").Value = 12098
$ws1.Range("F19").Value = 94

# ---- Sheet "全部类型" (same events, rows shifted down by one) ----
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F4").Value = 94
$ws4.Range("F6").Value = 45
$ws4.Range("F7").Value = 2676

$ws4.Range("C9").Value = "苏州·国风宠物-cosplay展（取消）"
$ws4.Range("G9").Value = "不可售"

$ws4.Range("F10").Value = 259
$ws4.Range("F11").Value = 116
$ws4.Range("F12").Value = 10043
$ws4.Range("F14").Value = 256
$ws4.Range("F15").Value = 8
$ws4.Range("F16").Value = 620
$ws4.Range("F17").Value = 11750
$ws4.Range("F18").Value = 12098
$ws4.Range("F20").Value = 94
